$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Updated time_taken values for rows 2-15 on the "data" sheet
$timestamps = @{
    2  = "2021-10-05 14:21:49.093275"
    3  = "2021-10-05 14:21:49.093285"
    4  = "2021-10-05 14:21:49.093288"
    5  = "2021-10-05 14:21:49.093291"
    6  = "2021-10-05 14:21:49.093295"
    7  = "2021-10-05 14:21:49.093297"
    8  = "2021-10-05 14:21:49.093300"
    9  = "2021-10-05 14:21:49.093303"
    10 = "2021-10-05 14:21:49.093306"
    11 = "2021-10-05 14:21:49.093309"
    12 = "2021-10-05 14:21:49.093311"
    13 = "2021-10-05 14:21:49.093314"
    14 = "2021-10-05 14:21:49.093317"
    15 = "2021-10-05 14:21:49.093319"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 6).Value = $timestamps[$row]
}

# Add a new "metadata" worksheet right after the "data" sheet
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1) - copy the header cell formatting from the "data" sheet
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

$ws.Range("B1:F1").Copy() | Out-Null
$meta.Range("B1:F1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Copy() | Out-Null
$meta.Range("G1").PasteSpecial(-4122) | Out-Null

# Data row (row 2)
$ws.Range("A2").Copy() | Out-Null
$meta.Range("A2").PasteSpecial(-4122) | Out-Null
$meta.Cells.Item(2, 1).Value = 0

$meta.Cells.Item(2, 2).Value = "Neuronal ceroid lipofuscinosis"
$meta.Cells.Item(2, 3).Value = 526

$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.23"
$meta.Range("D2").ClearFormats()

$meta.Cells.Item(2, 5).Value = "2021-03-17T18:27:10.592575Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:21:49.089740"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/526/?format=json"

$excel.CutCopyMode = 0
